# ---------------------------------------------------------------------------
# Applies two changes to TESTING ENTREAMIGOS.docx:
#   1. Highlights the "-SESION!!! ..." paragraph in magenta.
#   2. Removes the _GoBack bookmark from the end of the
#      "-Títulos a las páginas <h1>" paragraph and relocates it into the
#      middle of "una nueva" (between "un" and "a nueva") in the next
#      paragraph ("-Añadir una nueva organización como amigo no funciona
#      ..."), which otherwise keeps its original single-paragraph shape.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Highlight "-SESION!!!..." run in magenta (wdPink == 5 == magenta)
#    Setting HighlightColorIndex on the found (partial) range highlights
#    the whole enclosing run, same as in real Word.
# ---------------------------------------------------------------------------
$sesionRng = $d.Content
$sesionRng.Find.ClearFormatting()
$sesionFound = $sesionRng.Find.Execute("-SESION!!!")
if ($sesionFound) {
    $sesionRng.HighlightColorIndex = 5
}

# ---------------------------------------------------------------------------
# 2) Locate the "-Títulos a las páginas <h1>" paragraph and the following
#    "-Añadir una nueva organización ..." paragraph, then relocate the
#    _GoBack bookmark out of the former and into the middle of "una nueva"
#    in the latter.
# ---------------------------------------------------------------------------
$titlePara = $null
$addPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*T*tulos a las p*ginas*h1*") {
        $titlePara = $p
    }
    if ($t -like "*-A*adir*nueva*organizaci*n como amigo*funciona*") {
        $addPara = $p
    }
}

if (($titlePara -ne $null) -and ($addPara -ne $null)) {
    # Replace the whole two-paragraph span in one shot: this keeps the
    # "organización..." text glued to the same paragraph as "-Añadir un" +
    # bookmark + "a nueva " (InsertXML only merges cleanly at full
    # paragraph-range boundaries, so the "a nueva " / "organización" split
    # must NOT be expressed as a separate partial-range InsertXML call).
    $combined = $d.Range($titlePara.Range.Start, $addPara.Range.End)

    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
        + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
        + '<pkg:xmlData>' `
        + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
        + '<w:body>' `
        + '<w:p><w:r><w:t>-Títulos a las páginas &lt;h1&gt;</w:t></w:r></w:p>' `
        + '<w:p>' `
        + '<w:r><w:t xml:space="preserve">-Añadir </w:t></w:r>' `
        + '<w:r><w:t>un</w:t></w:r>' `
        + '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' `
        + '<w:bookmarkEnd w:id="0"/>' `
        + '<w:r><w:t xml:space="preserve">a nueva </w:t></w:r>' `
        + '<w:r><w:t>organización como amigo no funciona</w:t></w:r>' `
        + '<w:r><w:t xml:space="preserve"> a no ser que seas tu mismo (problema con nueva organización)</w:t></w:r>' `
        + '</w:p>' `
        + '</w:body></w:document>' `
        + '</pkg:xmlData></pkg:part></pkg:package>'

    $combined.InsertXML($pkg) | Out-Null
}
